# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (per-fund holdings detail) right before
# the "总计" (Total) summary sheet, and adds a corresponding summary row to
# the top of "总计".

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q1" sheet ---
# Clone the most recent per-quarter fund-detail sheet ("2021-Q4") so the new
# sheet inherits identical sheetPr/outline/page-setup formatting, inserting
# the clone right before "总计" (currently the last sheet).
#
# NOTE: once Copy()/Add() runs, any variable still pointing at the sheet
# that used to occupy that slot gets reseated onto the freshly-inserted
# sheet instead of staying bound to the original - so sheets are always
# re-fetched by name right before they're used, never held across an insert.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 already carries the correct per-column formatting inherited from
# the cloned template (A2 bold/bordered index style, B2:H2 plain). Fill in
# its values, then stamp identically-formatted copies down for rows 3-5
# before overwriting each row's values.
$newSheet.Range("B2").Value = "'011355"
$newSheet.Range("C2").Value = "华泰柏瑞港股通时代机遇混合型证券投资基金A"
$newSheet.Range("D2").Value = "'1.13"
$newSheet.Range("E2").Value = "'90.93"
$newSheet.Range("F2").Value = "'5.30"
$newSheet.Range("G2").Value = "'0.0599"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A2:H2").Copy($newSheet.Range("A3"))
$newSheet.Range("A2:H2").Copy($newSheet.Range("A4"))
$newSheet.Range("A2:H2").Copy($newSheet.Range("A5"))

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'003413"
$newSheet.Range("C3").Value = "华泰柏瑞新经济沪港深灵活配置混合"
$newSheet.Range("D3").Value = "'0.54"
$newSheet.Range("E3").Value = "'92.57"
$newSheet.Range("F3").Value = "'5.53"
$newSheet.Range("G3").Value = "'0.0299"
$newSheet.Range("H3").Value = 9

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'011356"
$newSheet.Range("C4").Value = "华泰柏瑞港股通时代机遇混合型证券投资基金C"
$newSheet.Range("D4").Value = "'0.40"
$newSheet.Range("E4").Value = "'90.93"
$newSheet.Range("F4").Value = "'5.30"
$newSheet.Range("G4").Value = "'0.0212"
$newSheet.Range("H4").Value = 9

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'002860"
$newSheet.Range("C5").Value = "前海开源沪港深新机遇灵活配置混合"
$newSheet.Range("D5").Value = "'0.01"
$newSheet.Range("E5").Value = "'83.26"
$newSheet.Range("F5").Value = "'6.90"
$newSheet.Range("G5").Value = "'0.0007"
$newSheet.Range("H5").Value = 3

# Drop the left-over "quote prefix" text-format marker picked up from
# forcing these numeric-looking values to stay text, so the cells end up
# plain (un-styled) like the rest of the data rows.
$newSheet.Range("B2:G5").ClearFormats()

# --- 2. Insert the new 2022-Q1 summary row at the top of "总计" ---
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The inserted row borrows row 3's (previously row 2's) index-column style.
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.11

# Renumber the 0-based index column for the rows pushed down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
